# Updates the cryptos list (Price / Volume(1h) columns, plus a few
# Coin/Link swaps) to match the latest scrape, as produced by the
# scheduled GitHub Actions job.
#
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
#
# Many "Price" values look like plain numbers (e.g. "11.52") which Excel
# would otherwise silently convert to numeric/date values on assignment.
# To preserve them as literal text (matching the original inlineStr
# cells), the cell's NumberFormat is forced to Text ("@") before writing
# those values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '22.395.88'
$ws.Cells.Item(2, 5).Value = '  +0.03%  '
$ws.Cells.Item(3, 4).Value = '1.573.21'
$ws.Cells.Item(3, 5).Value = '  +0.28%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '290.90'
$ws.Cells.Item(6, 5).Value = '  -0.08%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3767'
$ws.Cells.Item(7, 5).Value = '  +2.81%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '49.96'
$ws.Cells.Item(8, 5).Value = '  +1.21%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.3428'
$ws.Cells.Item(9, 5).Value = '  +1.42%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.07654'
$ws.Cells.Item(10, 5).Value = '  +0.86%  '
$ws.Cells.Item(11, 5).Value = '  -1.42%  '
$ws.Cells.Item(12, 5).Value = '  +0.09%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '21.26'
$ws.Cells.Item(13, 5).Value = '  +0.33%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.039'
$ws.Cells.Item(14, 5).Value = '  -0.16%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.953'
$ws.Cells.Item(15, 5).Value = '  +0.74%  '
$ws.Cells.Item(16, 4).Value = '1.573.44'
$ws.Cells.Item(16, 5).Value = '  +0.32%  '
$ws.Cells.Item(17, 5).Value = '  -0.34%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = '@'
$cell.Value = '90.22'
$ws.Cells.Item(18, 5).Value = '  +1.34%  '
$ws.Cells.Item(19, 5).Value = '  +0.10%  '
$ws.Cells.Item(20, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '16.85'
$ws.Cells.Item(21, 5).Value = '  +2.40%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.209'
$ws.Cells.Item(22, 5).Value = '  -0.35%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.03'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '
$ws.Cells.Item(24, 4).Value = '22.393.06'
$ws.Cells.Item(24, 5).Value = '  -0.06%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.402'
$ws.Cells.Item(25, 5).Value = '  +0.06%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.695'
$ws.Cells.Item(26, 5).Value = '  -10.16%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '20.23'
$ws.Cells.Item(27, 5).Value = '  +1.92%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = '@'
$cell.Value = '147.41'
$ws.Cells.Item(28, 5).Value = '  +1.55%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.038'
$ws.Cells.Item(29, 5).Value = '  +1.05%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '126.47'
$ws.Cells.Item(30, 5).Value = '  +0.90%  '
$ws.Cells.Item(31, 4).Value = '1.746.43'
$ws.Cells.Item(31, 5).Value = '  +0.09%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '6.185'
$ws.Cells.Item(32, 5).Value = '  -1.18%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.006'
$ws.Cells.Item(33, 5).Value = '  +0.82%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.9975'
$ws.Cells.Item(34, 5).Value = '  -4.06%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '10.01'
$ws.Cells.Item(35, 5).Value = '  -3.27%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.08567'
$ws.Cells.Item(36, 5).Value = '  +1.60%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.02555'
$ws.Cells.Item(37, 5).Value = '  +0.38%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.2320'
$ws.Cells.Item(38, 5).Value = '  +0.15%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.06586'
$ws.Cells.Item(39, 5).Value = '  +0.65%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.330'
$ws.Cells.Item(40, 5).Value = '  +6.17%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.447'
$ws.Cells.Item(41, 5).Value = '  -1.64%  '
$ws.Cells.Item(42, 2).Value = 'TheSandbox'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6430'
$ws.Cells.Item(42, 5).Value = '  +0.66%  '
$ws.Cells.Item(43, 2).Value = 'Aptos'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '11.52'
$ws.Cells.Item(43, 5).Value = '  -2.64%  '
$ws.Cells.Item(44, 5).Value = '  +0.09%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '14.05'
$ws.Cells.Item(45, 5).Value = '  -2.51%  '
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.6003'
$ws.Cells.Item(46, 5).Value = '  +0.11%  '
$ws.Cells.Item(47, 2).Value = 'PancakeSwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.784'
$ws.Cells.Item(47, 5).Value = '  +0.20%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.309'
$ws.Cells.Item(48, 5).Value = '  +7.47%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '2.091'
$ws.Cells.Item(49, 5).Value = '  -2.08%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '125.46'
$ws.Cells.Item(50, 5).Value = '  +1.84%  '
$ws.Cells.Item(51, 5).Value = '  +0.48%  '

Write-Host "Applied 99 cell updates"
